# Weekly update: a new price record (2022-07-21) is inserted as the new
# row 25 for "Vega Modelo de Temuco" / Maracuyá, pushing the existing
# rows 25-50 down to 26-51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 25, shifting rows 25:50 down to 26:51.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record's data.
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 44763
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100108
$ws.Range("H25").Value = "Tropicales y subtropicales"
$ws.Range("I25").Value = 100108003
$ws.Range("J25").Value = "Maracuyá"
$ws.Range("K25").Value = "Sin especificar"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 220
$ws.Range("N25").Value = 32000
$ws.Range("O25").Value = 34000
$ws.Range("P25").Value = 33091
$ws.Range("Q25").Value = "$/caja 18 kilos"
$ws.Range("R25").Value = "Región de Arica y Parinacota"
$ws.Range("S25").Value = 1838
$ws.Range("T25").Value = 18
